# Applies the "Lagt till 106 och 108 ..." progress-log update to Blad1.
#
# Summary of data changes (per the commit's xml diff):
#  - Row 6  (Distance sampling-analys på ripdata): status ej påbörjat -> klar,
#    new comment in C6 about qq-plots.
#  - Row 9  (lyaktivitet reproduktion): status påbörjat -> klar, comment
#    replaced with a short note that Lars corrected the files.
#  - Row 10 (GIS-data lyornas avstånd trädgräns): comment updated to mention
#    Blankan and Norr Vaktklumpen.
#  - Row 11 (GIS-data lyornas avstånd vatten...): comment updated likewise.
#  - Row 12 (GIS-data area myrar...): comment updated likewise.
#  - Row 15 (GIS-data lyornas höjd över havet): new notes added in C15/D15.
#  - Selection moves from C8 to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Distance sampling-analys på ripdata -> klar, plus new comment ---
$klarColor = $ws.Range("B2").Font.Color
$ws.Range("B6").Value = "klar"
$ws.Range("B6").Font.Color = $klarColor
$ws.Range("C6").Value = "är qq-plot etc tillräckligt bra för ripskitar?"

# --- Row 9: lyaktivitet reproduktion -> klar, comment replaced ---
$ws.Range("B9").Value = "klar"
$ws.Range("B9").Font.Color = $klarColor
$ws.Range("C9").Value = "Lars har rättat mina filer."

# --- Row 10: mention Blankan och Norr Vaktklumpen in the comment ---
$ws.Range("C10").Value = "Lagt till Blankan och Norr Vaktklumpen! tog bort små trädsamlingar som inte satt ihop med skog som går ner i låglandet."

# --- Row 11: mention Blankan och Norr vaktklumpen in the comment ---
$ws.Range("C11").Value = "Blankan och Norr vaktklumpen klar. Klart för närmsta vattenkälla för alla lyor. Tog andel vatten inom 1,5 km radie istället för närmsta avstånd till större vatten."

# --- Row 12: mention Blankan och Norr vaktklumpen in the comment ---
$ws.Range("C12").Value = "Blankan och Norr vaktklumpen klar. Myrar (vadare och lämmel). Tar andel myr inom en radie på 1,5 km (samma avstånd som Gallant et al 2014 gjorde för sina parametrar). KOLLA SÅ ATT ALLT ÄR I EPSG:3006 Sweref! Tar cirkel, inte triangel."

# --- Row 15: new notes in C15 / D15 ---
$ws.Range("C15").Value = "Blankan och Norr vaktklumpen klar"
$ws.Range("D15").Value = "pen klar"

# --- Move the active selection to C10 (was C8) ---
$ws.Range("C10").Select()
